# Update the "relatorio_neomater_COMPLETO" report values to reflect the
# re-generated report after adding the "competencia" (period) and "ano"
# (year) filters when building the report (fixes a bug where the
# presentation values in Excel were not found/filled correctly).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - PACOTE PRÉ-OPERATÓRIO PEDIÁTRICO OTORRINO
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0

# Row 3 - PACOTE PRÉ-OPERATÓRIO PEDIÁTRICO CIRURGIA GERAL
$ws.Range("B3").Value = 11
$ws.Range("E3").Value = 7
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 7
$ws.Range("I3").Value = 5
$ws.Range("K3").Value = 12
$ws.Range("L3").Value = 6
$ws.Range("M3").Value = 9

# Row 5 - ADENOIDECTOMIA PEDIÁTRICO
$ws.Range("E5").Value = 0
$ws.Range("L5").Value = 0

# Row 6 - AMIGDALECTOMIA- PEDIATRICO
$ws.Range("C6").Value = 0

# Row 7 - AMIGDALECTOMIA COM ADENOIDECTOMIA - PEDIATRICO
$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 0
$ws.Range("G7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0

# Row 10 - HERNIOPLASTIA INGUINAL (BILATERAL) - PEDIATRICO
$ws.Range("B10").Value = 4
$ws.Range("E10").Value = 1
$ws.Range("G10").Value = 1
$ws.Range("K10").Value = 4
$ws.Range("L10").Value = 2
$ws.Range("M10").Value = 1

# Row 11 - HERNIOPLASTIA UMBILICAL - PEDIATRICO
$ws.Range("B11").Value = 7
$ws.Range("E11").Value = 4
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 3
$ws.Range("I11").Value = 5
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 2
$ws.Range("M11").Value = 2

# Row 12 - ORQUIDOPEXIA BILATERAL - PEDIATRICO
$ws.Range("B12").Value = 1
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 0

# Row 14 - CORRECAO DE HIPOSPADIA (1º TEMPO) - PEDIATRICO
$ws.Range("B14").Value = 1
$ws.Range("G14").Value = 1
$ws.Range("K14").Value = 2

# Row 16 - POSTECTOMIA - PEDIATRICO
$ws.Range("B16").Value = 7
$ws.Range("E16").Value = 4
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 5
$ws.Range("K16").Value = 8
$ws.Range("L16").Value = 5
$ws.Range("M16").Value = 8

# Row 17 - TOTAL
$ws.Range("B17").Value = 31
$ws.Range("C17").Value = 0
$ws.Range("D17").Value = 0
$ws.Range("E17").Value = 16
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = 17
$ws.Range("I17").Value = 11
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 30
$ws.Range("M17").Value = 20
